$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 284, shifting existing rows 284-299 down to 285-300
$ws.Rows("284:284").Insert()

# Populate the newly inserted row 284 with its data
$ws.Range("A284").Value = 9
$ws.Range("B284").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C284").Value = "Metropolitana"
$ws.Range("D284").Value = 44585
$ws.Range("E284").Value = 13
$ws.Range("F284").Value = 100112032
$ws.Range("G284").Value = "Zapallo italiano"
$ws.Range("H284").Value = "Sin especificar"
$ws.Range("I284").Value = "Primera"
$ws.Range("J284").Value = 61
$ws.Range("K284").Value = 13000
$ws.Range("L284").Value = 14000
$ws.Range("M284").Value = 13508
$ws.Range("N284").Value = "$/caja 50 unidades"
$ws.Range("O284").Value = "Región de O'Higgins"
$ws.Range("P284").Value = 270
$ws.Range("Q284").Value = 50
$ws.Range("R284").Value = "Hortaliza"
